$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data (Session 11 with second strategy)
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = 20

# Move selection to A13, mirroring the user clicking the next empty cell
$ws.Range("A13").Select()
